# The source data added two new weekly observations (rows 412-413) for
# "Albahaca" at Mercado Mayorista Lo Valledor de Santiago, pushing the
# previously-existing rows 412-479 down to 414-481 (dimension grows from
# A1:R479 to A1:R481).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at 412 and 413, shifting everything below
# (old rows 412..479) down to 414..481, carrying formatting along (this
# also extends the used range / dimension automatically).
$ws.Rows.Item(412).EntireRow.Insert()
$ws.Rows.Item(413).EntireRow.Insert()

# New row 412: Primera quality observation dated 2022-05-13 (serial 44694)
$ws.Cells.Item(412, 1).Value = 6
$ws.Cells.Item(412, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(412, 3).Value = "Metropolitana"
$ws.Cells.Item(412, 4).Value = 44694
$ws.Cells.Item(412, 5).Value = 13
$ws.Cells.Item(412, 6).Value = 100112052
$ws.Cells.Item(412, 7).Value = "Albahaca"
$ws.Cells.Item(412, 8).Value = "Sin especificar"
$ws.Cells.Item(412, 9).Value = "Primera"
$ws.Cells.Item(412, 10).Value = 80
$ws.Cells.Item(412, 11).Value = 2500
$ws.Cells.Item(412, 12).Value = 3000
$ws.Cells.Item(412, 13).Value = 2688
$ws.Cells.Item(412, 14).Value = "$/docena de matas"
$ws.Cells.Item(412, 15).Value = "Región Metropolitana"
$ws.Cells.Item(412, 16).Value = 448
$ws.Cells.Item(412, 17).Value = 6
$ws.Cells.Item(412, 18).Value = "Hortaliza"

# New row 413: Segunda quality observation, same date
$ws.Cells.Item(413, 1).Value = 6
$ws.Cells.Item(413, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(413, 3).Value = "Metropolitana"
$ws.Cells.Item(413, 4).Value = 44694
$ws.Cells.Item(413, 5).Value = 13
$ws.Cells.Item(413, 6).Value = 100112052
$ws.Cells.Item(413, 7).Value = "Albahaca"
$ws.Cells.Item(413, 8).Value = "Sin especificar"
$ws.Cells.Item(413, 9).Value = "Segunda"
$ws.Cells.Item(413, 10).Value = 50
$ws.Cells.Item(413, 11).Value = 2000
$ws.Cells.Item(413, 12).Value = 2000
$ws.Cells.Item(413, 13).Value = 2000
$ws.Cells.Item(413, 14).Value = "$/docena de matas"
$ws.Cells.Item(413, 15).Value = "Región Metropolitana"
$ws.Cells.Item(413, 16).Value = 333
$ws.Cells.Item(413, 17).Value = 6
$ws.Cells.Item(413, 18).Value = "Hortaliza"

Write-Host "Inserted rows 412-413; new dimension ref is now A1:R481"
